$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, column B)
$ws.Range("B8").Value = "2025-08-20T17:48:34+01:00"

# Set the Description value (row 12, column B), which was previously empty
$ws.Range("B12").Value = "Extended value set for sleep quality assessment"
